$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the "A" column labels for the existing 2023 Duke standard-curve rows
# (standard_curve_2023 -> duke_standard_curve_2023 .. duke_standard_curve_2027)
$ws.Range("A2").Value = "duke_standard_curve_2023"
$ws.Range("A3").Value = "duke_standard_curve_2024"
$ws.Range("A4").Value = "duke_standard_curve_2025"
$ws.Range("A5").Value = "duke_standard_curve_2026"
$ws.Range("A6").Value = "duke_standard_curve_2027"

# Append the new Duke August 2023 standard curve block (rows 12-16)
$ws.Range("A12").Value = "duke_august_standard_curve_2023"
$ws.Range("A13").Value = "duke_august_standard_curve_2023"
$ws.Range("A14").Value = "duke_august_standard_curve_2023"
$ws.Range("A15").Value = "duke_august_standard_curve_2023"
$ws.Range("A16").Value = "duke_august_standard_curve_2023"

$ws.Range("C12").Value = "Duke"
$ws.Range("C13").Value = "Duke"
$ws.Range("C14").Value = "Duke"
$ws.Range("C15").Value = "Duke"
$ws.Range("C16").Value = "Duke"

$ws.Range("D12").Value = "BernhardtTriology"
$ws.Range("D13").Value = "BernhardtTriology"
$ws.Range("D14").Value = "BernhardtTriology"
$ws.Range("D15").Value = "BernhardtTriology"
$ws.Range("D16").Value = "BernhardtTriology"

$ws.Range("E12").Value = "low_zero"
$ws.Range("E13").Value = "high_full"
$ws.Range("E14").Value = "high_zero"
$ws.Range("E15").Value = "cutoff"
$ws.Range("E16").Value = "above_detection"

$ws.Range("F12").Value = 0.2351
$ws.Range("F13").Value = 0.437
$ws.Range("F14").Value = 0.3325
$ws.Range("F15").Value = 3977.04
$ws.Range("F16").Value = 10661.8

$ws.Range("G12").Value = 0
$ws.Range("G13").Value = -877
$ws.Range("G14").Value = 0

$ws.Range("H12").Value = 2023
$ws.Range("H13").Value = 2023
$ws.Range("H14").Value = 2023
$ws.Range("H15").Value = 2023
$ws.Range("H16").Value = 2023

$ws.Range("I12").Value = 0.995
$ws.Range("I13").Value = 0.96
$ws.Range("I14").Value = 0.989

$ws.Range("J13").Value = "y/x = rfu/ug per L, high > 1000 ug per L or > 3977rfu"
$ws.Range("J14").Value = "y/x = rfu/ug per L, high > 1000 ug per L or > 3977rfu"
$ws.Range("J12").Value = "y/x = rfu/ug per L, low < 1000 ug per L or <3977rfu"
$ws.Range("J15").Value = "max values in 1000 ug per L"
$ws.Range("J16").Value = "max value in 4000 ug/L"

$ws.Range("B12").Value = "DukeAug2023"
$ws.Range("B13").Value = "DukeAug2023"
$ws.Range("B14").Value = "DukeAug2023"
$ws.Range("B15").Value = "DukeAug2023"
$ws.Range("B16").Value = "DukeAug2023"

# Update the active-cell selection to match the post-edit workbook view
$ws.Range("B18").Select()
